$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.3464964993005633
$ws.Range("C2").Value = 0.3375848360084654
$ws.Range("D2").Value = 0.1529057820181812
$ws.Range("E2").Value = 71517.89157740913
$ws.Range("G2").Value = 71518.72856452645

$ws.Range("B3").Value = 0.0000001063418937352623
$ws.Range("C3").Value = 0.05231270169004087
$ws.Range("D3").Value = 0.7127328510149897
$ws.Range("E3").Value = 198602002.3250627
$ws.Range("G3").Value = 198602003.0901083
